# Applies todo-list updates: renamed ColorPalette library references & new
# backlog items (refactor: update to match renamed library ColorPalette to Colors)
$wb = $excel.ActiveWorkbook

# --- Sheet "Active" (sheet1): rewrite task rows 2-26 ---
$wsActive = $wb.Worksheets.Item("Active")

# Format the "Created" column as Text first so date-like strings (e.g. "8/9/2018")
# are stored as literal text instead of being auto-converted to date serials.
$wsActive.Range("E2:E26").NumberFormat = "@"

$wsActive.Range("A2").Value = 23
$wsActive.Range("B2").Value = "what do these lines do?`nApplication.EnableVisualStyles();`nApplication.SetCompatibleTextRenderingDefault(false);"
$wsActive.Range("C2").Value = "Todo"
$wsActive.Range("D2").Value = "Task"
$wsActive.Range("E2").Value = "8/9/2018"

$wsActive.Range("A3").Value = 36
$wsActive.Range("B3").Value = "what to rename ColorPalette library to so it does not conflict with object ColorPalette?"
$wsActive.Range("C3").Value = "Todo"
$wsActive.Range("D3").Value = "Task"
$wsActive.Range("E3").Value = "8/12/2018"

$wsActive.Range("A4").Value = 41
$wsActive.Range("B4").Value = "bug: on some colors (oranges esp.) adjusting the saturation gives a too bright color - keep it in the gray range"
$wsActive.Range("C4").Value = "Todo"
$wsActive.Range("D4").Value = "Task"
$wsActive.Range("E4").Value = "8/15/2018"

$wsActive.Range("A5").Value = 38
$wsActive.Range("B5").Value = "support applying color while zoomed in or out"
$wsActive.Range("C5").Value = "Todo"
$wsActive.Range("D5").Value = "Task"
$wsActive.Range("E5").Value = "8/12/2018"

$wsActive.Range("A6").Value = 33
$wsActive.Range("B6").Value = "undo, redo coloring a section on the image"
$wsActive.Range("C6").Value = "Todo"
$wsActive.Range("D6").Value = "Task"
$wsActive.Range("E6").Value = "8/11/2018"

$wsActive.Range("A7").Value = 37
$wsActive.Range("B7").Value = "apply color over color without changing underlying values`n- change section back to grayscale`n- then to the new color"
$wsActive.Range("C7").Value = "Todo"
$wsActive.Range("D7").Value = "Task"
$wsActive.Range("E7").Value = "8/12/2018"

$wsActive.Range("A8").Value = 44
$wsActive.Range("B8").Value = "move coloring operations into another thread`n- queue incoming commands`n- handle one at a time in another thread, updating display between each one"
$wsActive.Range("C8").Value = "Todo"
$wsActive.Range("D8").Value = "Task"
$wsActive.Range("E8").Value = "8/15/2018"

$wsActive.Range("A9").Value = 45
$wsActive.Range("B9").Value = "display closable modal ""Please Wait"" while coloring image`n- with ""Cancel Color"" option that will stop the thread and cancel the changes"
$wsActive.Range("C9").Value = "Todo"
$wsActive.Range("D9").Value = "Task"
$wsActive.Range("E9").Value = "8/15/2018"

$wsActive.Range("A10").Value = 42
$wsActive.Range("B10").Value = "bug: expanded palette covers part of picturebox"
$wsActive.Range("C10").Value = "Todo"
$wsActive.Range("D10").Value = "Task"
$wsActive.Range("E10").Value = "8/15/2018"

$wsActive.Range("A11").Value = 35
$wsActive.Range("B11").Value = "set and check tolerance for ""black"" and ""white""`n- ""blacks"" will be left untouched`n- ""whites"" will be treated as pure white, which will in effect alter them to white"
$wsActive.Range("C11").Value = "Todo"
$wsActive.Range("D11").Value = "Task"
$wsActive.Range("E11").Value = "8/12/2018"

$wsActive.Range("A12").Value = 13
$wsActive.Range("B12").Value = "open edit palette mode`n- add swatches`n- remove swatches`n- change color of swatch`n- reorder swatches`n- undo, redo until pane is closed`n- save changes before leaving pane `n(save as text file - or, what do other programs use?)`n(no duplicate colors allowed)"
$wsActive.Range("C12").Value = "Todo"
$wsActive.Range("D12").Value = "Task"
$wsActive.Range("E12").Value = "8/9/2018"

$wsActive.Range("A13").Value = 34
$wsActive.Range("B13").Value = "when zooming, if a scroll bar is all the way to min or max, keep it there"
$wsActive.Range("C13").Value = "Todo"
$wsActive.Range("D13").Value = "Task"
$wsActive.Range("E13").Value = "8/11/2018"

$wsActive.Range("A14").Value = 21
$wsActive.Range("B14").Value = "when resizing windows, default behavior is to keep the same section of image in the viewing pane`n- so expanding window would in effect zoom in"
$wsActive.Range("C14").Value = "Todo"
$wsActive.Range("D14").Value = "Task"
$wsActive.Range("E14").Value = "8/9/2018"

$wsActive.Range("A15").Value = 24
$wsActive.Range("B15").Value = "remember windows size from last closing`n- full screen vs not`n- default not-full-screen size`nopen with this size"
$wsActive.Range("C15").Value = "Todo"
$wsActive.Range("D15").Value = "Task"
$wsActive.Range("E15").Value = "8/9/2018"

$wsActive.Range("A16").Value = 27
$wsActive.Range("B16").Value = "include support contact information`n-wohaste email`n-paint landing page on website`n-github page"
$wsActive.Range("C16").Value = "Todo"
$wsActive.Range("D16").Value = "Task"
$wsActive.Range("E16").Value = "8/9/2018"

$wsActive.Range("A17").Value = 28
$wsActive.Range("B17").Value = "include donation information`n-patreon"
$wsActive.Range("C17").Value = "Todo"
$wsActive.Range("D17").Value = "Task"
$wsActive.Range("E17").Value = "8/9/2018"

$wsActive.Range("A18").Value = 19
$wsActive.Range("B18").Value = "design an icon`napply to windows, desktop icon, and uninstall icon"
$wsActive.Range("C18").Value = "Todo"
$wsActive.Range("D18").Value = "Task"
$wsActive.Range("E18").Value = "8/9/2018"

$wsActive.Range("A19").Value = 22
$wsActive.Range("B19").Value = "check through the program Properties >> Assembly Info again, make sure all is accurate"
$wsActive.Range("C19").Value = "Todo"
$wsActive.Range("D19").Value = "Task"
$wsActive.Range("E19").Value = "8/9/2018"

$wsActive.Range("A20").Value = 16
$wsActive.Range("B20").Value = "that should be it for minimum viable`nsave as Version 1"
$wsActive.Range("C20").Value = "Todo"
$wsActive.Range("D20").Value = "Task"
$wsActive.Range("E20").Value = "8/9/2018"

$wsActive.Range("A21").Value = 17
$wsActive.Range("B21").Value = "build installer for Version 1 and save it in a separate folder to keep"
$wsActive.Range("C21").Value = "Todo"
$wsActive.Range("D21").Value = "Task"
$wsActive.Range("E21").Value = "8/9/2018"

$wsActive.Range("A22").Value = 18
$wsActive.Range("B22").Value = "update website with project, landing page, and links"
$wsActive.Range("C22").Value = "Todo"
$wsActive.Range("D22").Value = "Task"
$wsActive.Range("E22").Value = "8/9/2018"

$wsActive.Range("A23").Value = 25
$wsActive.Range("B23").Value = "EVERYTHING BELOW HERE IS VERSION 2"
$wsActive.Range("C23").Value = "Todo"
$wsActive.Range("D23").Value = "Task"
$wsActive.Range("E23").Value = "8/9/2018"

$wsActive.Range("A24").Value = 26
$wsActive.Range("B24").Value = "remember last used directory (save or open) and default to there in file dialogs"
$wsActive.Range("C24").Value = "Todo"
$wsActive.Range("D24").Value = "Task"
$wsActive.Range("E24").Value = "8/9/2018"

$wsActive.Range("A25").Value = 43
$wsActive.Range("B25").Value = "on deep zoom, when image is getting too big, switch to zooming on just a segment of the masterImage`n- this will complicate scrollbars and color placement"
$wsActive.Range("C25").Value = "Todo"
$wsActive.Range("D25").Value = "Task"
$wsActive.Range("E25").Value = "8/15/2018"

$wsActive.Range("A26").Value = 29
$wsActive.Range("B26").Value = "how to programs auto-update?"
$wsActive.Range("C26").Value = "Todo"
$wsActive.Range("D26").Value = "Task"
$wsActive.Range("E26").Value = "8/9/2018"

# --- Sheet "Inactive" (sheet2): rewrite task rows 2-14 ---
$wsInactive = $wb.Worksheets.Item("Inactive")

$wsInactive.Range("E2:E14").NumberFormat = "@"
$wsInactive.Range("F2:F14").NumberFormat = "@"

$wsInactive.Range("A2").Value = 11
$wsInactive.Range("B2").Value = "save changes to image`n- bitmap`n- png`n- jpg"
$wsInactive.Range("C2").Value = "Done"
$wsInactive.Range("D2").Value = "Task"
$wsInactive.Range("E2").Value = "8/9/2018"
$wsInactive.Range("F2").Value = "8/15/2018"

$wsInactive.Range("A3").Value = 40
$wsInactive.Range("B3").Value = "zoom in needs to keep pixels clear instead of letting it blur together"
$wsInactive.Range("C3").Value = "Done"
$wsInactive.Range("D3").Value = "Task"
$wsInactive.Range("E3").Value = "8/13/2018"
$wsInactive.Range("F3").Value = "8/15/2018"

$wsInactive.Range("A4").Value = 39
$wsInactive.Range("B4").Value = "bug: it isn't actually keeping the grayscale`n- wow, spent days debugging and it was just a test-line outside the area I was looking at"
$wsInactive.Range("C4").Value = "Done"
$wsInactive.Range("D4").Value = "Task"
$wsInactive.Range("E4").Value = "8/12/2018"
$wsInactive.Range("F4").Value = "8/15/2018"

$wsInactive.Range("A5").Value = 32
$wsInactive.Range("B5").Value = "fill in a section of color on the image"
$wsInactive.Range("C5").Value = "Done"
$wsInactive.Range("D5").Value = "Task"
$wsInactive.Range("E5").Value = "8/11/2018"
$wsInactive.Range("F5").Value = "8/12/2018"

$wsInactive.Range("A6").Value = 31
$wsInactive.Range("B6").Value = "select a palette color"
$wsInactive.Range("C6").Value = "Done"
$wsInactive.Range("D6").Value = "Task"
$wsInactive.Range("E6").Value = "8/11/2018"
$wsInactive.Range("F6").Value = "8/11/2018"

$wsInactive.Range("A7").Value = 20
$wsInactive.Range("B7").Value = "program preference setting: how wide the palette area is`n- let user drag and drop divider to change palette width`n- minimum is 3 swatches, max is maybe 12`n- save setting and reuse on next progam open`n(maybe instead of drag-n-drop the border, there are little < > arrow buttons that will expand/contract space one swatch at a time)"
$wsInactive.Range("C7").Value = "Done"
$wsInactive.Range("D7").Value = "Task"
$wsInactive.Range("E7").Value = "8/9/2018"
$wsInactive.Range("F7").Value = "8/11/2018"

$wsInactive.Range("A8").Value = 9
$wsInactive.Range("B8").Value = "display a default palette along the side of the window`n- give it a vertical scroll when needed"
$wsInactive.Range("C8").Value = "Done"
$wsInactive.Range("D8").Value = "Task"
$wsInactive.Range("E8").Value = "8/9/2018"
$wsInactive.Range("F8").Value = "8/11/2018"

$wsInactive.Range("A9").Value = 8
$wsInactive.Range("B9").Value = "build 1 to 3 default palettes`n- downloaded some Photoshop palettes"
$wsInactive.Range("C9").Value = "Done"
$wsInactive.Range("D9").Value = "Task"
$wsInactive.Range("E9").Value = "8/9/2018"
$wsInactive.Range("F9").Value = "8/11/2018"

$wsInactive.Range("A10").Value = 4
$wsInactive.Range("B10").Value = "scroll horizontal and vertical when zoomed in"
$wsInactive.Range("C10").Value = "Done"
$wsInactive.Range("D10").Value = "Task"
$wsInactive.Range("E10").Value = "8/9/2018"
$wsInactive.Range("F10").Value = "8/10/2018"

$wsInactive.Range("A11").Value = 3
$wsInactive.Range("B11").Value = "zoom in and out of image"
$wsInactive.Range("C11").Value = "Done"
$wsInactive.Range("D11").Value = "Task"
$wsInactive.Range("E11").Value = "8/9/2018"
$wsInactive.Range("F11").Value = "8/9/2018"

$wsInactive.Range("A12").Value = 30
$wsInactive.Range("B12").Value = "resize image to fit window as it resizes"
$wsInactive.Range("C12").Value = "Done"
$wsInactive.Range("D12").Value = "Task"
$wsInactive.Range("E12").Value = "8/9/2018"
$wsInactive.Range("F12").Value = "8/9/2018"

$wsInactive.Range("A13").Value = 2
$wsInactive.Range("B13").Value = "open an image and display it in a window at default size"
$wsInactive.Range("C13").Value = "Done"
$wsInactive.Range("D13").Value = "Task"
$wsInactive.Range("E13").Value = "8/9/2018"
$wsInactive.Range("F13").Value = "8/9/2018"

$wsInactive.Range("A14").Value = 1
$wsInactive.Range("B14").Value = "new window's console app`n- version 0 until minimum viable is complete"
$wsInactive.Range("C14").Value = "Done"
$wsInactive.Range("D14").Value = "Task"
$wsInactive.Range("E14").Value = "8/9/2018"
$wsInactive.Range("F14").Value = "8/9/2018"

# --- Sheet "Config" (sheet3): bump Max Id to the new highest task id (45) ---
$wsConfig = $wb.Worksheets.Item("Config")
$wsConfig.Range("F2").Value = 45
